$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2380.0715
$ws.Range("I43").Value = 4320.2
$ws.Range("J43").Value = 1302.2222
$ws.Range("K43").Value = 4320.2
$ws.Range("L43").Value = 1302.2222
$ws.Range("M43").Value = -4251.2
$ws.Range("N43").Value = -1440.2222

$ws.Range("H70").Value = 33335184
$ws.Range("I70").Value = 625.5
$ws.Range("J70").Value = 45456844
$ws.Range("K70").Value = 1876.5
$ws.Range("L70").Value = 136370532
$ws.Range("M70").Value = -1606.5
$ws.Range("N70").Value = -136371072

$ws.Range("H73").Value = 33335184
$ws.Range("I73").Value = 625.5
$ws.Range("J73").Value = 45456844
$ws.Range("K73").Value = 1876.5
$ws.Range("L73").Value = 136370532
$ws.Range("M73").Value = -940.5
$ws.Range("N73").Value = -136372404

$ws.Range("H74").Value = 5662265.5
$ws.Range("I74").Value = 7278770
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 7278770
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -7277834
$ws.Range("N74").Value = -6372

$ws.Range("H77").Value = 5662265.5
$ws.Range("I77").Value = 7278770
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 36393850
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -36389170
$ws.Range("N77").Value = -31860

$ws.Range("H98").Value = 1706.55
$ws.Range("I98").Value = 1851.5625
$ws.Range("J98").Value = 1126.5
$ws.Range("K98").Value = 1851.5625
$ws.Range("L98").Value = 1126.5
$ws.Range("M98").Value = -353.5625
$ws.Range("N98").Value = -4122.5

$ws.Range("H122").Value = 1706.55
$ws.Range("I122").Value = 1851.5625
$ws.Range("J122").Value = 1126.5
$ws.Range("K122").Value = 5554.6875
$ws.Range("L122").Value = 3379.5
$ws.Range("M122").Value = -3104.6875
$ws.Range("N122").Value = -8279.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2871.6365
$ws.Range("I88").Value = 2400.25
$ws.Range("J88").Value = 3141
$ws.Range("K88").Value = 2400.25
$ws.Range("L88").Value = 3141
$ws.Range("M88").Value = -1994.25
$ws.Range("N88").Value = -3953

$ws.Range("H91").Value = 2871.6365
$ws.Range("I91").Value = 2400.25
$ws.Range("J91").Value = 3141
$ws.Range("K91").Value = 2400.25
$ws.Range("L91").Value = 3141
$ws.Range("M91").Value = -996.25
$ws.Range("N91").Value = -5949

$ws.Range("H97").Value = 1191.9
$ws.Range("I97").Value = 1202.5
$ws.Range("J97").Value = 1149.5
$ws.Range("K97").Value = 1202.5
$ws.Range("L97").Value = 1149.5
$ws.Range("M97").Value = -706.5
$ws.Range("N97").Value = -2141.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 774.29034
$ws.Range("I94").Value = 791.76666
$ws.Range("J94").Value = 250
$ws.Range("K94").Value = 791.76666
$ws.Range("L94").Value = 250
$ws.Range("M94").Value = -340.76666
$ws.Range("N94").Value = -1152

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37708.586
$ws.Range("I31").Value = 3630.5293
$ws.Range("J31").Value = 85985.836
$ws.Range("K31").Value = 3630.5293
$ws.Range("L31").Value = 85985.836
$ws.Range("M31").Value = -3335.5293
$ws.Range("N31").Value = -86575.836

$ws.Range("H34").Value = 37708.586
$ws.Range("I34").Value = 3630.5293
$ws.Range("J34").Value = 85985.836
$ws.Range("K34").Value = 3630.5293
$ws.Range("L34").Value = 85985.836
$ws.Range("M34").Value = -3428.5293
$ws.Range("N34").Value = -86389.836

$ws.Range("H62").Value = 4682.273
$ws.Range("J62").Value = 4277.778
$ws.Range("L62").Value = 4277.778
$ws.Range("N62").Value = -5525.778

$ws.Range("H65").Value = 4682.273
$ws.Range("J65").Value = 4277.778
$ws.Range("L65").Value = 21388.89
$ws.Range("N65").Value = -27628.89

$ws.Range("H107").Value = 458.22223
$ws.Range("I107").Value = 521.8333
$ws.Range("J107").Value = 426.41666
$ws.Range("K107").Value = 521.8333
$ws.Range("L107").Value = 426.41666
$ws.Range("M107").Value = 1398.1667
$ws.Range("N107").Value = -4266.41666

$ws.Range("H132").Value = 1771.8
$ws.Range("I132").Value = 1262.4
$ws.Range("J132").Value = 2281.2
$ws.Range("K132").Value = 3787.2
$ws.Range("L132").Value = 6843.599999999999
$ws.Range("M132").Value = -1257.2
$ws.Range("N132").Value = -11903.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4167202.8
$ws.Range("J68").Value = 567.0714
$ws.Range("L68").Value = 1701.2142
$ws.Range("N68").Value = -3323.2142

$ws.Range("H71").Value = 4167202.8
$ws.Range("J71").Value = 567.0714
$ws.Range("L71").Value = 5103.6426
$ws.Range("N71").Value = -13215.6426

$ws.Range("H107").Value = 339065.1
$ws.Range("I107").Value = 1120
$ws.Range("J107").Value = 458339.8
$ws.Range("K107").Value = 3360
$ws.Range("L107").Value = 1375019.4
$ws.Range("M107").Value = -1440
$ws.Range("N107").Value = -1378859.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1512.0769
$ws.Range("I97").Value = 1808.8889
$ws.Range("J97").Value = 844.25
$ws.Range("K97").Value = 1808.8889
$ws.Range("L97").Value = 844.25
$ws.Range("M97").Value = -1312.8889
$ws.Range("N97").Value = -1836.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2872.8262
$ws.Range("I93").Value = 2491.6667
$ws.Range("J93").Value = 3117.8572
$ws.Range("K93").Value = 2491.6667
$ws.Range("L93").Value = 3117.8572
$ws.Range("M93").Value = -1243.6667
$ws.Range("N93").Value = -5613.8572

$ws.Range("H132").Value = 2126.2727
$ws.Range("I132").Value = 1971.4546
$ws.Range("J132").Value = 2590.7273
$ws.Range("K132").Value = 5914.3638
$ws.Range("L132").Value = 7772.1819
$ws.Range("M132").Value = -3384.3638
$ws.Range("N132").Value = -12832.1819

$ws.Range("H136").Value = 2960.0393
$ws.Range("I136").Value = 1815.45
$ws.Range("K136").Value = 5446.35
$ws.Range("M136").Value = -2896.35

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1605.6666
$ws.Range("I96").Value = 1231.6666
$ws.Range("J96").Value = 2166.6667
$ws.Range("K96").Value = 1231.6666
$ws.Range("L96").Value = 2166.6667
$ws.Range("M96").Value = 141.3334
$ws.Range("N96").Value = -4912.6667

$ws.Range("H132").Value = 1029.4857
$ws.Range("I132").Value = 878.5517
$ws.Range("K132").Value = 2635.6551
$ws.Range("M132").Value = -105.6550999999999

$ws.Range("H136").Value = 1160.9111
$ws.Range("I136").Value = 1069.1578
$ws.Range("K136").Value = 3207.4734
$ws.Range("M136").Value = -657.4733999999999
